$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (attendee interest count) for three events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 285
$wsExhibit.Range("F4").Value = 958
$wsExhibit.Range("F6").Value = 55

# Sheet "全部类型" - same three events repeated, update the same field
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 285
$wsAll.Range("F5").Value = 958
$wsAll.Range("F7").Value = 55
